$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("oc-files")

# Insert a new column before column C
$ws.Columns("C").Insert()

# Fill data rows with "Optical Clearing" first (so this string is registered first)
$ws.Range("C2:C90").Value = "Optical Clearing"

# Set header for new column
$ws.Range("C1").Value = "Image Category"

# Update selection to match target state
$ws.Range("C1").Select()
